# "added logo on IdCardGen"
# Replace the per-row photo URL text shown in the PhotoURL column with a
# generic default image placeholder ("default.jpg") while leaving the
# existing hyperlinks (which still point at the original photo URLs)
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "default.jpg"
$ws.Range("E3").Value = "default.jpg"
$ws.Range("E4").Value = "default.jpg"

# Leave the final selection on the last edited cell.
$ws.Range("E4").Select() | Out-Null
